# fixed bug in the fedrollover
# Append the next 7 "rollover" placeholder rows (Jan-Jul 2009) that were
# missing from the auction-dates sheet, following the same pattern as the
# existing trailing placeholder rows (cusip/total/offering/percents = 0,
# only the rollover "date" column populated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("20090100", "20090200", "20090300", "20090400", "20090500", "20090600", "20090700")

$startRow = 636
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 0
    # Leading apostrophe forces the rollover-date code to be stored as text
    # (matches the "date" column type used by every other row in the sheet)
    # instead of being auto-coerced to a number; re-applying the "Normal"
    # style afterwards drops the transient quote-prefix formatting so the
    # cell ends up unstyled, same as its neighbours.
    $ws.Cells.Item($r, 2).Value = "'" + $dates[$i]
    $ws.Cells.Item($r, 2).Style = "Normal"
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
}
